$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the two missing collision-matrix entries (playerBox x opponentBox)
# using the same formatting as the other "f" cells in the sheet.
$ws.Range("H5").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("F5").Value = "f"

$ws.Range("C6").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Value = "f"

$excel.CutCopyMode = $false

# Update the active selection on the sheet.
$ws.Range("G7").Select()
